# feat: add 2022-Q1 data
#
# The workbook has sheets: 2020-Q4, 2021-Q1, 2021-Q2, 2021-Q3, 2021-Q4, 总计
# We need to:
#   1. Insert a new "2022-Q1" sheet (fund holdings detail) right before "总计"
#   2. Update "总计" (summary) with a new row for 2022-Q1 at the top, shifting
#      all other rows down by one.
#
# Implementation notes:
#  - We duplicate the existing "总计" sheet (Copy) so the new summary sheet
#    inherits the same sheet-level formatting (sheetPr/pageMargins/etc).
#  - The ORIGINAL "总计" worksheet object is renamed to "2022-Q1" and its
#    data is cleared & replaced with the fund-holding detail rows. This keeps
#    its original sheetId, so the new sheet ordering/IDs match what a normal
#    "insert before" would produce.
#  - The DUPLICATE sheet (freshly minted, gets the next sheetId) is renamed
#    back to "总计" and repopulated with the refreshed summary table.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")

# Duplicate "总计" right after itself -- the copy will become the new "总计"
# summary sheet, while we repurpose the original sheet object as "2022-Q1".
$total.Copy([System.Reflection.Missing]::Value, $total)
$newTotal = $wb.Worksheets.Item($total.Index + 1)

# --- Build the "2022-Q1" sheet (fund holdings detail) out of the original
# --- "总计" sheet object ---
$q1 = $total
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$headerSrc = $wb.Worksheets.Item("2021-Q4")

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
$headerSrc.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q1Data = @(
    @("011052", "鹏华弘裕一年持有期混合A", "2.92", "24.56", "0.66", "0.0193", 9),
    @("011053", "鹏华弘裕一年持有期混合C", "0.20", "24.56", "0.66", "0.0013", 9)
)

$r = 2
foreach ($row in $q1Data) {
    $q1.Range("A$r").Value = ($r - 2)
    $q1.Range("B$r").Value = "'" + $row[0]
    $q1.Range("B$r").Style = "Normal"
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = "'" + $row[2]
    $q1.Range("D$r").Style = "Normal"
    $q1.Range("E$r").Value = "'" + $row[3]
    $q1.Range("E$r").Style = "Normal"
    $q1.Range("F$r").Value = "'" + $row[4]
    $q1.Range("F$r").Style = "Normal"
    $q1.Range("G$r").Value = "'" + $row[5]
    $q1.Range("G$r").Style = "Normal"
    $q1.Range("H$r").Value = $row[6]
    $r = $r + 1
}

$headerSrc.Range("A2:A3").Copy()
$q1.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Build the refreshed "总计" summary sheet out of the duplicated sheet ---
$newTotal.Name = "总计"
$newTotal.Cells.Clear()

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"
$headerSrc.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalData = @(
    @("2022-Q1", 2, 0.02),
    @("2021-Q4", 19, 4.08),
    @("2021-Q3", 7, 1.32),
    @("2021-Q2", 15, 2.47),
    @("2021-Q1", 10, 0.5),
    @("2020-Q4", 4, 0.66)
)

$r = 2
foreach ($row in $totalData) {
    $newTotal.Range("A$r").Value = ($r - 2)
    $newTotal.Range("B$r").Value = $row[0]
    $newTotal.Range("C$r").Value = $row[1]
    $newTotal.Range("D$r").Value = $row[2]
    $r = $r + 1
}

$headerSrc.Range("A2:A7").Copy()
$newTotal.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the original active sheet/selection state (unchanged by this edit).
$wb.Worksheets.Item("2020-Q4").Activate()
